$wb = $excel.ActiveWorkbook

# This script re-applies updated currentAveragePrice / LevePrice / LeveProfit
# figures (columns H-N) across several Leve-profit sheets, as pulled from the
# latest market-board data by the scheduled runner.

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2896.5862
$ws.Cells.Item(15, 9).Value = 2896.5862
$ws.Cells.Item(15, 11).Value = 8689.758600000001
$ws.Cells.Item(15, 13).Value = -8520.758600000001

$ws.Cells.Item(99, 8).Value = 874.8125
$ws.Cells.Item(99, 9).Value = 673.7778
$ws.Cells.Item(99, 10).Value = 1133.2858
$ws.Cells.Item(99, 11).Value = 2021.3334
$ws.Cells.Item(99, 12).Value = 3399.8574
$ws.Cells.Item(99, 13).Value = -523.3334
$ws.Cells.Item(99, 14).Value = -6395.857400000001

$ws.Cells.Item(115, 8).Value = 7214.2173
$ws.Cells.Item(115, 10).Value = 9504.706
$ws.Cells.Item(115, 12).Value = 28514.118
$ws.Cells.Item(115, 14).Value = -31648.118

$ws.Cells.Item(137, 8).Value = 13598598
$ws.Cells.Item(137, 9).Value = 915.76666
$ws.Cells.Item(137, 10).Value = 39094252
$ws.Cells.Item(137, 11).Value = 2747.29998
$ws.Cells.Item(137, 12).Value = 117282756
$ws.Cells.Item(137, 13).Value = -197.2999799999998
$ws.Cells.Item(137, 14).Value = -117287856

$ws.Cells.Item(138, 8).Value = 2000.3536
$ws.Cells.Item(138, 9).Value = 1352.2106
$ws.Cells.Item(138, 10).Value = 3478.12
$ws.Cells.Item(138, 11).Value = 4056.6318
$ws.Cells.Item(138, 12).Value = 10434.36
$ws.Cells.Item(138, 13).Value = 1083.3682
$ws.Cells.Item(138, 14).Value = -20714.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 27879926
$ws.Cells.Item(74, 9).Value = 27027620
$ws.Cells.Item(74, 10).Value = 29631888
$ws.Cells.Item(74, 11).Value = 27027620
$ws.Cells.Item(74, 12).Value = 29631888
$ws.Cells.Item(74, 13).Value = -27026746
$ws.Cells.Item(74, 14).Value = -29633636

$ws.Cells.Item(77, 8).Value = 27879926
$ws.Cells.Item(77, 9).Value = 27027620
$ws.Cells.Item(77, 10).Value = 29631888
$ws.Cells.Item(77, 11).Value = 135138100
$ws.Cells.Item(77, 12).Value = 148159440
$ws.Cells.Item(77, 13).Value = -135133732
$ws.Cells.Item(77, 14).Value = -148168176

$ws.Cells.Item(110, 8).Value = 2309.3572
$ws.Cells.Item(110, 9).Value = 1284.4286
$ws.Cells.Item(110, 10).Value = 3334.2856
$ws.Cells.Item(110, 11).Value = 1284.4286
$ws.Cells.Item(110, 12).Value = 3334.2856
$ws.Cells.Item(110, 13).Value = 760.5714
$ws.Cells.Item(110, 14).Value = -7424.2856

$ws.Cells.Item(132, 8).Value = 10514421
$ws.Cells.Item(132, 9).Value = 11115345
$ws.Cells.Item(132, 10).Value = 7939035.5
$ws.Cells.Item(132, 11).Value = 33346035
$ws.Cells.Item(132, 12).Value = 23817106.5
$ws.Cells.Item(132, 13).Value = -33343505
$ws.Cells.Item(132, 14).Value = -23822166.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1789754.2
$ws.Cells.Item(31, 9).Value = 1840.409
$ws.Cells.Item(31, 10).Value = 4815454.5
$ws.Cells.Item(31, 11).Value = 1840.409
$ws.Cells.Item(31, 12).Value = 4815454.5
$ws.Cells.Item(31, 13).Value = -1545.409
$ws.Cells.Item(31, 14).Value = -4816044.5

$ws.Cells.Item(34, 8).Value = 1789754.2
$ws.Cells.Item(34, 9).Value = 1840.409
$ws.Cells.Item(34, 10).Value = 4815454.5
$ws.Cells.Item(34, 11).Value = 1840.409
$ws.Cells.Item(34, 12).Value = 4815454.5
$ws.Cells.Item(34, 13).Value = -1638.409
$ws.Cells.Item(34, 14).Value = -4815858.5

$ws.Cells.Item(58, 8).Value = 1517618
$ws.Cells.Item(58, 9).Value = 2977.913
$ws.Cells.Item(58, 10).Value = 6494292.5
$ws.Cells.Item(58, 11).Value = 2977.913
$ws.Cells.Item(58, 12).Value = 6494292.5
$ws.Cells.Item(58, 13).Value = -2774.913
$ws.Cells.Item(58, 14).Value = -6494698.5

$ws.Cells.Item(134, 8).Value = 956479.9399999999
$ws.Cells.Item(134, 9).Value = 3912.6572
$ws.Cells.Item(134, 10).Value = 5719316.5
$ws.Cells.Item(134, 11).Value = 11737.9716
$ws.Cells.Item(134, 12).Value = 17157949.5
$ws.Cells.Item(134, 13).Value = -9202.971600000001
$ws.Cells.Item(134, 14).Value = -17163019.5

$ws.Cells.Item(136, 8).Value = 1517618
$ws.Cells.Item(136, 9).Value = 2977.913
$ws.Cells.Item(136, 10).Value = 6494292.5
$ws.Cells.Item(136, 11).Value = 8933.739
$ws.Cells.Item(136, 12).Value = 19482877.5
$ws.Cells.Item(136, 13).Value = -6383.739
$ws.Cells.Item(136, 14).Value = -19487977.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 6486759.5
$ws.Cells.Item(5, 9).Value = 15385524
$ws.Cells.Item(5, 10).Value = 2778940.8
$ws.Cells.Item(5, 11).Value = 46156572
$ws.Cells.Item(5, 12).Value = 8336822.399999999
$ws.Cells.Item(5, 13).Value = -46156460
$ws.Cells.Item(5, 14).Value = -8337046.399999999

$ws.Cells.Item(92, 8).Value = 4878824.5
$ws.Cells.Item(92, 9).Value = 280
$ws.Cells.Item(92, 10).Value = 6098461
$ws.Cells.Item(92, 11).Value = 840
$ws.Cells.Item(92, 12).Value = 18295383
$ws.Cells.Item(92, 13).Value = 408
$ws.Cells.Item(92, 14).Value = -18297879

$ws.Cells.Item(107, 8).Value = 480.1579
$ws.Cells.Item(107, 9).Value = 244.28572
$ws.Cells.Item(107, 10).Value = 617.75
$ws.Cells.Item(107, 11).Value = 732.85716
$ws.Cells.Item(107, 12).Value = 1853.25
$ws.Cells.Item(107, 13).Value = 1187.14284
$ws.Cells.Item(107, 14).Value = -5693.25

$ws.Cells.Item(129, 8).Value = 53573836
$ws.Cells.Item(129, 9).Value = 104167750
$ws.Cells.Item(129, 10).Value = 10207629
$ws.Cells.Item(129, 11).Value = 312503250
$ws.Cells.Item(129, 12).Value = 30622887
$ws.Cells.Item(129, 13).Value = -312498250
$ws.Cells.Item(129, 14).Value = -30632887

$ws.Cells.Item(132, 8).Value = 2674.7368
$ws.Cells.Item(132, 10).Value = 2664
$ws.Cells.Item(132, 12).Value = 23976
$ws.Cells.Item(132, 14).Value = -29036

$ws.Cells.Item(135, 8).Value = 6486759.5
$ws.Cells.Item(135, 9).Value = 15385524
$ws.Cells.Item(135, 10).Value = 2778940.8
$ws.Cells.Item(135, 11).Value = 138469716
$ws.Cells.Item(135, 12).Value = 25010467.2
$ws.Cells.Item(135, 13).Value = -138467181
$ws.Cells.Item(135, 14).Value = -25015537.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(141, 8).Value = 44851.145
$ws.Cells.Item(141, 10).Value = 44851.145
$ws.Cells.Item(141, 12).Value = 44851.145
$ws.Cells.Item(141, 14).Value = -55211.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4085977
$ws.Cells.Item(132, 9).Value = 4766056.5
$ws.Cells.Item(132, 11).Value = 14298169.5
$ws.Cells.Item(132, 13).Value = -14295639.5

$ws.Cells.Item(136, 8).Value = 3368653
$ws.Cells.Item(136, 9).Value = 5051896
$ws.Cells.Item(136, 10).Value = 2166.3635
$ws.Cells.Item(136, 11).Value = 15155688
$ws.Cells.Item(136, 12).Value = 6499.0905
$ws.Cells.Item(136, 13).Value = -15153138
$ws.Cells.Item(136, 14).Value = -11599.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 8476.757
$ws.Cells.Item(107, 9).Value = 9321.182000000001
$ws.Cells.Item(107, 10).Value = 7238.2666
$ws.Cells.Item(107, 11).Value = 27963.546
$ws.Cells.Item(107, 12).Value = 21714.7998
$ws.Cells.Item(107, 13).Value = -26043.546
$ws.Cells.Item(107, 14).Value = -25554.7998

$ws.Cells.Item(132, 8).Value = 620417.9399999999
$ws.Cells.Item(132, 9).Value = 1724.881
$ws.Cells.Item(132, 10).Value = 4332576
$ws.Cells.Item(132, 11).Value = 5174.643
$ws.Cells.Item(132, 12).Value = 12997728
$ws.Cells.Item(132, 13).Value = -2644.643
$ws.Cells.Item(132, 14).Value = -13002788

$ws.Cells.Item(140, 8).Value = 58186.668
$ws.Cells.Item(140, 10).Value = 58186.668
$ws.Cells.Item(140, 12).Value = 58186.668
$ws.Cells.Item(140, 14).Value = -68546.66800000001

